$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "mngr577004"
$ws.Range("B2").Value = "rYtUneb"
